$wb = $excel.ActiveWorkbook

$wsChangePin = $wb.Worksheets.Item("ChangePin")
$wsVerifyPin = $wb.Worksheets.Item("VerifyPin")

# --- ChangePin sheet ---
# A9 had the string "null"; clear it to an empty cell (keeps style)
$wsChangePin.Range("A9").ClearContents()

# Move the active selection to A9
$wsChangePin.Activate()
$wsChangePin.Range("A9").Select()

# --- VerifyPin sheet ---
# Update the id/value rows from 1 -> 155
$wsVerifyPin.Range("A2").Value = 155
$wsVerifyPin.Range("A7").Value = 155
$wsVerifyPin.Range("A8").Value = 155
$wsVerifyPin.Range("A9").Value = 155
$wsVerifyPin.Range("A10").Value = 155
$wsVerifyPin.Range("A11").Value = 155
$wsVerifyPin.Range("A12").Value = 155
$wsVerifyPin.Range("A13").Value = 155

# A6 and B12 had the string "null"; clear them to empty cells
$wsVerifyPin.Range("A6").ClearContents()
$wsVerifyPin.Range("B12").ClearContents()

# B13 value change 123457 -> 123456
$wsVerifyPin.Range("B13").Value = 123456

# Move the active selection to A14
$wsVerifyPin.Activate()
$wsVerifyPin.Range("A14").Select()
